{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// worksheet table with its new value. Old values are unique in the\n// document, so a body.search() per pair unambiguously locates the\n// single run to update.\nconst replacements = [\n  [\"572\u00d74=\", \"295\u00d72=\"],\n  [\"537\u00d75=\", \"113\u00d76=\"],\n  [\"416\u00d76=\", \"834\u00d73=\"],\n  [\"726\u00d75=\", \"647\u00d73=\"],\n  [\"262\u00d77=\", \"157\u00d76=\"],\n  [\"867\u00d77=\", \"634\u00d77=\"],\n  [\"378\u00d73=\", \"841\u00d77=\"],\n  [\"447\u00d72=\", \"814\u00d74=\"],\n  [\"281\u00d74=\", \"209\u00d74=\"],\n  [\"574\u00d74=\", \"481\u00d77=\"],\n  [\"858\u00d77=\", \"439\u00d75=\"],\n  [\"808\u00d78=\", \"359\u00d79=\"],\n  [\"624\u00d74=\", \"226\u00d79=\"],\n  [\"637\u00d77=\", \"983\u00d74=\"],\n  [\"945\u00d74=\", \"586\u00d77=\"],\n  [\"766\u00d76=\", \"683\u00d72=\"],\n  [\"367\u00d77=\", \"583\u00d78=\"],\n  [\"706\u00d78=\", \"440\u00d76=\"],\n  [\"776\u00d73=\", \"897\u00d77=\"],\n  [\"550\u00d75=\", \"683\u00d72=\"],\n  [\"358\u00d78=\", \"516\u00d78=\"],\n  [\"702\u00d76=\", \"677\u00d76=\"],\n  [\"182\u00d74=\", \"649\u00d73=\"],\n  [\"754\u00d75=\", \"798\u00d73=\"],\n  [\"644\u00d79=\", \"913\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# worksheet table with its new value. Old values are unique in the\n# document, so a simple Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"572\u00d74=\"; New = \"295\u00d72=\" },\n    @{ Old = \"537\u00d75=\"; New = \"113\u00d76=\" },\n    @{ Old = \"416\u00d76=\"; New = \"834\u00d73=\" },\n    @{ Old = \"726\u00d75=\"; New = \"647\u00d73=\" },\n    @{ Old = \"262\u00d77=\"; New = \"157\u00d76=\" },\n    @{ Old = \"867\u00d77=\"; New = \"634\u00d77=\" },\n    @{ Old = \"378\u00d73=\"; New = \"841\u00d77=\" },\n    @{ Old = \"447\u00d72=\"; New = \"814\u00d74=\" },\n    @{ Old = \"281\u00d74=\"; New = \"209\u00d74=\" },\n    @{ Old = \"574\u00d74=\"; New = \"481\u00d77=\" },\n    @{ Old = \"858\u00d77=\"; New = \"439\u00d75=\" },\n    @{ Old = \"808\u00d78=\"; New = \"359\u00d79=\" },\n    @{ Old = \"624\u00d74=\"; New = \"226\u00d79=\" },\n    @{ Old = \"637\u00d77=\"; New = \"983\u00d74=\" },\n    @{ Old = \"945\u00d74=\"; New = \"586\u00d77=\" },\n    @{ Old = \"766\u00d76=\"; New = \"683\u00d72=\" },\n    @{ Old = \"367\u00d77=\"; New = \"583\u00d78=\" },\n    @{ Old = \"706\u00d78=\"; New = \"440\u00d76=\" },\n    @{ Old = \"776\u00d73=\"; New = \"897\u00d77=\" },\n    @{ Old = \"550\u00d75=\"; New = \"683\u00d72=\" },\n    @{ Old = \"358\u00d78=\"; New = \"516\u00d78=\" },\n    @{ Old = \"702\u00d76=\"; New = \"677\u00d76=\" },\n    @{ Old = \"182\u00d74=\"; New = \"649\u00d73=\" },\n    @{ Old = \"754\u00d75=\"; New = \"798\u00d73=\" },\n    @{ Old = \"644\u00d79=\"; New = \"913\u00d73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $r = $d.Content\n    $r.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
